# Add two new header columns (I0, IF) mirroring the existing header style,
# and their corresponding data-row values, extending the sheet from H to J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header cell formatting (bold font, thin border, centered/top
# aligned) from the existing last header cell (H1) onto the two new header
# cells so they pick up the same style index used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for row 2
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7
